# Update localization status report for archive.
#
# 1. Status value "Ready for handoff" -> "In Translation" on all three
#    sheets (Overview!E2/F2, zh-cn!C2, de-de!C2 - they all point at the
#    same shared string).
# 2. Narrow the "Status"/"zh-cn"/"de-de" status columns from
#    17.2159881591797 to 13.4101845877511 characters wide:
#      - Overview: columns E and F
#      - zh-cn:    column C
#      - de-de:    column C

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# The OOXML <col width="..."> is derived from the COM ColumnWidth
# (character units) by snapping to whole pixels, so the literal target
# width of 13.4101845877511 can't be typed in directly - 12.5 is the
# ColumnWidth that lands closest to it after Excel's pixel rounding.
$newWidth  = 12.5

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E:E").ColumnWidth = $newWidth
$overview.Range("F:F").ColumnWidth = $newWidth

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C:C").ColumnWidth = $newWidth

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C:C").ColumnWidth = $newWidth
